# Updated remaining queries for C3DC
#
# Applies the JOIN-clause "id" -> "study_id" / "participant_id" rename
# across every TabQuery / StatQuery cell on the sheet, plus the
# sheetView/topLeftCell scroll position and the column C width tweak
# captured in the diff.
#
# NOTE on style: this COM host's PowerShell function-parameter binding for
# *named* parameters (`-Text $x`) does not propagate the argument into the
# function scope correctly, so helper functions below are always invoked
# positionally, and every COM property write goes through a plain local
# variable first (never `<expr>.Prop = Func args` in one statement) to
# dodge that issue entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Update-QueryText {
    param([string]$Text)

    $result = $Text
    $result = $result.Replace('std.id = prt."study.id"', 'std.study_id = prt."study.study_id"')
    $result = $result.Replace('prt.id = dgn."participant.id"', 'prt.participant_id = dgn."participant.participant_id"')
    $result = $result.Replace('prt.id = trt."participant.id"', 'prt.participant_id = trt."participant.participant_id"')
    $result = $result.Replace('prt.id = trr."participant.id"', 'prt.participant_id = trr."participant.participant_id"')
    $result = $result.Replace('prt.id = srv."participant.id"', 'prt.participant_id = srv."participant.participant_id"')
    $result = $result.Replace('std.id = rfs."study.id"', 'std.study_id = rfs."study.study_id"')

    return $result
}

# StudiesTab (B2/C2), ParticipantsTab (B3), DiagnosisTab (B4),
# TreatmentTab (B5), TreatmentRespTab (B6), SurvivalTab (B7)
$cellsToUpdate = @("B2", "C2", "B3", "B4", "B5", "B6", "B7")

foreach ($addr in $cellsToUpdate) {
    $cell = $ws.Range($addr)
    $oldText = $cell.Value2
    $newText = Update-QueryText $oldText
    $cell.Value2 = $newText
}

# sheetView topLeftCell scrolled back up to A2 (was A6)
$activeWindow = $excel.ActiveWindow
$activeWindow.ScrollRow = 2
$activeWindow.ScrollColumn = 1

# Column C widened from 60.83 to 69.83 characters and no longer "best fit"
$colC = $ws.Columns("C")
$colC.ColumnWidth = 69
